# Weekly price update: a new observation is inserted as row 447 (pushing every
# subsequent historical row down by one), matching the "Fruta / hortaliza,
# semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 447; Excel shifts rows
# 447:501 down to 448:502 and extends the used range to A1:R502 for us.
$ws.Rows.Item(447).Insert()

# Populate the newly-inserted row 447 with the new weekly data point.
$ws.Range("A447").Value = 8
$ws.Range("B447").Value = "Terminal La Palmera de La Serena"
$ws.Range("C447").Value = "Coquimbo"
$ws.Range("D447").Value = 45131
$ws.Range("E447").Value = 4
$ws.Range("F447").Value = 100112003
$ws.Range("G447").Value = "Ajo"
$ws.Range("H447").Value = "Chino"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 400
$ws.Range("K447").Value = 19500
$ws.Range("L447").Value = 20000
$ws.Range("M447").Value = 19750
$ws.Range("N447").Value = "$/caja 10 kilos"
$ws.Range("O447").Value = "China"
$ws.Range("P447").Value = 1975
$ws.Range("Q447").Value = 10
$ws.Range("R447").Value = "Hortaliza"
